# Apply "Natmi following Dr Hou advice" update:
# Adds ECs and FAPs as sending clusters (5 target clusters each,
# including two new target clusters M1 and M2) for the Wnt2-Fzd2 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Wnt2", "Fzd2", "ECs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.006017198313602724, 0.006017198313602724, 3, 1, 0.08013433333333332, 0.240403, 0.008122697422248188, 0.008122697422248188, 0.0008576243467777776, 0.007718619120999999, 0.00004887588123105699, 0.00004887588123105699),
    @("ECs", "Wnt2", "Fzd2", "FAPs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.006017198313602724, 0.006017198313602724, 3, 1, 7.647982, 22.943946, 0.7752263117781459, 0.775226311778146, 0.08185125269133332, 0.7366612742219999, 0.00466469045589192, 0.00466469045589192),
    @("ECs", "Wnt2", "Fzd2", "M1", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.006017198313602724, 0.006017198313602724, 1, 0.3333333333333333, 0.004539333333333334, 0.013618, 0.0004601227667548901, 0.00046012276675489, 0.00004858145844444444, 0.000437233126, 0.000002768649936167744, 0.000002768649936167744),
    @("ECs", "Wnt2", "Fzd2", "M2", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.006017198313602724, 0.006017198313602724, 1, 0.3333333333333333, 0.009266333333333333, 0.027799, 0.0009392680858436766, 0.0009392680858436766, 0.0000991713881111111, 0.000892542493, 0.00000565176234215943, 0.00000565176234215943),
    @("ECs", "Wnt2", "Fzd2", "sCs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.006017198313602724, 0.006017198313602724, 3, 1, 2.123561, 6.370683, 0.2152515999470071, 0.2152515999470071, 0.02272705767566667, 0.204543519081, 0.001295211564201419, 0.001295211564201419),
    @("FAPs", "Wnt2", "Fzd2", "ECs", 3, 1, 1.767921666666667, 5.303765, 0.9939828016863973, 0.9939828016863973, 3, 1, 0.08013433333333332, 0.240403, 0.008122697422248188, 0.008122697422248188, 0.1416712241438889, 1.275041017295, 0.00807382154101713, 0.00807382154101713),
    @("FAPs", "Wnt2", "Fzd2", "FAPs", 3, 1, 1.767921666666667, 5.303765, 0.9939828016863973, 0.9939828016863973, 3, 1, 7.647982, 22.943946, 0.7752263117781459, 0.775226311778146, 13.52103308407667, 121.68929775669, 0.770561621322254, 0.7705616213222541),
    @("FAPs", "Wnt2", "Fzd2", "M1", 3, 1, 1.767921666666667, 5.303765, 0.9939828016863973, 0.9939828016863973, 1, 0.3333333333333333, 0.004539333333333334, 0.013618, 0.0004601227667548901, 0.00046012276675489, 0.008025185752222223, 0.07222667177, 0.0004573541168187224, 0.0004573541168187223),
    @("FAPs", "Wnt2", "Fzd2", "M2", 3, 1, 1.767921666666667, 5.303765, 0.9939828016863973, 0.9939828016863973, 1, 0.3333333333333333, 0.009266333333333333, 0.027799, 0.0009392680858436766, 0.0009392680858436766, 0.01638215147055556, 0.147439363235, 0.0009336163235015172, 0.0009336163235015172),
    @("FAPs", "Wnt2", "Fzd2", "sCs", 3, 1, 1.767921666666667, 5.303765, 0.9939828016863973, 0.9939828016863973, 3, 1, 2.123561, 6.370683, 0.2152515999470071, 0.2152515999470071, 3.754289502388334, 33.788605521495, 0.2139563883828057, 0.2139563883828057)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $excelRow = $i + 2
    $rowdata = $rows[$i]
    for ($j = 0; $j -lt $rowdata.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $rowdata[$j]
    }
}
